$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Overwrite row 16 (currently "DANNA MARCELA HERNANDEZ DIAZ") with the data
# that previously lived in row 19 ("ANA ELVIRA CHICO PADILLA"), keeping row
# 16's own formatting. Then delete the old rows 17-19 entirely (which
# removes the now-duplicated/obsolete employee rows and shifts everything
# below up by 3).

$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1143381697"
$ws.Range("D16").Value = "ANA ELVIRA CHICO PADILLA"
$ws.Range("E16").Value = "2504"
$ws.Range("F16").Value = 11388
$ws.Range("G16").Value = 1423500

# Delete old rows 17, 18, 19 (NOHELY, DERLIS, and the now-duplicate ANA ELVIRA row)
$ws.Range("A17:A19").EntireRow.Delete()

# Update the summary figures
$ws.Range("E11").Value = 11388
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1

# Column D's width is auto-fit to content; with the longer names gone, let
# Excel recompute the best-fit width for the remaining (shorter) text.
$ws.Columns.Item(4).AutoFit()
